$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header strings for columns L:S (row 1)
$headers = @(
    "hzj-混合调节_20170516_152754_ASIC_EEG",
    "hzj-混合调节_20170518_134207_ASIC_EEG",
    "hzj-混合调节_20170519_135415_ASIC_EEG",
    "zyx-混合调节_20170516_111228_ASIC_EEG",
    "zyx-混合调节_20170517_110944_ASIC_EEG",
    "zyx-混合调节_20170518_112337_ASIC_EEG",
    "zyx-混合调节_20170519_124954_ASIC_EEG",
    "zyx-混合调节_20170522_111557_ASIC_EEG"
)

# New values for row 2 (columns L:S)
$row2 = @(
    0.94791666666666663,
    0.91752577319587625,
    0.90636704119850187,
    0.91987179487179482,
    0.91346153846153844,
    0.89967637540453071,
    0.9258064516129032,
    0.91747572815533984
)

# New values for row 3 (columns L:S)
$row3 = @(
    0.97802197802197799,
    0.93548387096774199,
    0.93040293040293043,
    0.93928571428571428,
    0.91061452513966479,
    0.89347079037800692,
    0.91275167785234901,
    0.90909090909090906
)

$cols = @("L", "M", "N", "O", "P", "Q", "R", "S")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $col = $cols[$i]
    $ws.Range("$col`1").Value = $headers[$i]
    $ws.Range("$col`2").Value = $row2[$i]
    $ws.Range("$col`3").Value = $row3[$i]
}

$ws.Range("A1:S3").Select()
